$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.308902025222778
$ws.Range("B1").Value = 1.867549657821655
$ws.Range("C1").Value = 1.74571681022644
$ws.Range("D1").Value = 1.693768501281738
$ws.Range("E1").Value = 1.166829586029053
